$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(390).Insert()

$ws.Range("A390").Value = 3
$ws.Range("B390").Value = "Femacal de La Calera"
$ws.Range("C390").Value = "Coquimbo"
$ws.Range("D390").Value = 44783
$ws.Range("E390").Value = 5
$ws.Range("F390").Value = 100112031
$ws.Range("G390").Value = "Poroto verde"
$ws.Range("H390").Value = "Magnum"
$ws.Range("I390").Value = "Primera"
$ws.Range("J390").Value = 78
$ws.Range("K390").Value = 32000
$ws.Range("L390").Value = 33000
$ws.Range("M390").Value = 32513
$ws.Range("N390").Value = '$/malla 25 kilos'
$ws.Range("O390").Value = "Región de Arica y Parinacota"
$ws.Range("P390").Value = 1301
$ws.Range("Q390").Value = 25
$ws.Range("R390").Value = "Hortaliza"
